$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 and J1
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy style from H1 (bold/border/center) onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# Data values for columns I (I0) and J (IF), rows 2-59
$data = @(
    @(2, 9, 9),
    @(3, 9, 9),
    @(4, 9, 9),
    @(5, 8, 9),
    @(6, 8, 9),
    @(7, 9, 9),
    @(8, 7, 7),
    @(9, 9, 9),
    @(10, 8, 8),
    @(11, 8, 9),
    @(12, 8, 9),
    @(13, 8, 9),
    @(14, 7, 7),
    @(15, 8, 8),
    @(16, 7, 8),
    @(17, 8, 9),
    @(18, 8, 9),
    @(19, 8, 8),
    @(20, 9, 9),
    @(21, 7, 7),
    @(22, 9, 9),
    @(23, 9, 9),
    @(24, 10, 10),
    @(25, 8, 8),
    @(26, 7, 7),
    @(27, 7, 8),
    @(28, 8, 8),
    @(29, 8, 8),
    @(30, 7, 8),
    @(31, 7, 8),
    @(32, 5, 5),
    @(33, 7, 7),
    @(34, 10, 10),
    @(35, 6, 6),
    @(36, 5, 7),
    @(37, 6, 7),
    @(38, 5, 5),
    @(39, 7, 7),
    @(40, 6, 6),
    @(41, 7, 7),
    @(42, 6, 7),
    @(43, 6, 7),
    @(44, 9, 9),
    @(45, 7, 7),
    @(46, 5, 5),
    @(47, 5, 6),
    @(48, 7, 7),
    @(49, 9, 9),
    @(50, 11, 12),
    @(51, 7, 7),
    @(52, 7, 7),
    @(53, 8, 9),
    @(54, 4, 5),
    @(55, 6, 6),
    @(56, 7, 7),
    @(57, 7, 8),
    @(58, 6, 7),
    @(59, 8, 8)
)

foreach ($entry in $data) {
    $r = $entry[0]
    $iVal = $entry[1]
    $jVal = $entry[2]
    $ws.Cells.Item($r, 9).Value = $iVal
    $ws.Cells.Item($r, 10).Value = $jVal
}
